$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, pushing existing rows 49:70 down to 50:71
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new data record
$ws.Cells.Item(49, 1).Value = 3
$ws.Cells.Item(49, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(49, 3).Value = "Coquimbo"
$ws.Cells.Item(49, 4).Value = 44784
$ws.Cells.Item(49, 5).Value = 5
$ws.Cells.Item(49, 6).Value = 100112035
$ws.Cells.Item(49, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 105
$ws.Cells.Item(49, 11).Value = 14000
$ws.Cells.Item(49, 12).Value = 15000
$ws.Cells.Item(49, 13).Value = 14476
$ws.Cells.Item(49, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(49, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(49, 16).Value = 965
$ws.Cells.Item(49, 17).Value = 15
$ws.Cells.Item(49, 18).Value = "Hortaliza"
